$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 39
$ws.Range("D2").Value = 1

$ws.Range("B5").Value = 0.975
$ws.Range("D5").Value = 0.025
